# Update the training plan schedule: shift Workout ID numbering and
# refresh the WT (workout type) / duration values for the affected rows,
# then apply a 2-decimal number format to the microcycle/weekly summary
# Miles / Duration / Est. Mileage columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workout ID (column A) for rows 15-126 is always (row - 14) ---
for ($r = 15; $r -le 126; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 14
}

# --- WT (column F) / duration (column H) overrides for specific rows ---
$overrides = @(
    @{ Row = 16; F = "T3"; H = 45 },
    @{ Row = 17; F = $null; H = 45 },
    @{ Row = 19; F = "T20"; H = 60 },
    @{ Row = 20; F = $null; H = 45 },
    @{ Row = 21; F = $null; H = 45 },
    @{ Row = 22; F = "E4x10s"; H = $null },
    @{ Row = 25; F = $null; H = 48 },
    @{ Row = 26; F = $null; H = 48 },
    @{ Row = 28; F = $null; H = 48 },
    @{ Row = 29; F = $null; H = 48 },
    @{ Row = 31; F = $null; H = 48 },
    @{ Row = 35; F = $null; H = 51 },
    @{ Row = 36; F = $null; H = 51 },
    @{ Row = 38; F = $null; H = 51 },
    @{ Row = 39; F = $null; H = 51 },
    @{ Row = 41; F = $null; H = 51 },
    @{ Row = 45; F = $null; H = 54 },
    @{ Row = 46; F = $null; H = 54 },
    @{ Row = 48; F = $null; H = 54 },
    @{ Row = 49; F = $null; H = 54 },
    @{ Row = 51; F = $null; H = 54 },
    @{ Row = 56; F = "T3"; H = 45 },
    @{ Row = 57; F = $null; H = 58 },
    @{ Row = 59; F = "T20"; H = 60 },
    @{ Row = 60; F = $null; H = 58 },
    @{ Row = 61; F = $null; H = 58 },
    @{ Row = 62; F = "E4x10s"; H = 45 },
    @{ Row = 65; F = $null; H = 61 },
    @{ Row = 66; F = $null; H = 61 },
    @{ Row = 68; F = $null; H = 61 },
    @{ Row = 69; F = $null; H = 61 },
    @{ Row = 71; F = $null; H = 61 },
    @{ Row = 75; F = $null; H = 64 },
    @{ Row = 76; F = $null; H = 64 },
    @{ Row = 78; F = $null; H = 64 },
    @{ Row = 79; F = $null; H = 64 },
    @{ Row = 81; F = $null; H = 64 },
    @{ Row = 85; F = $null; H = 67 },
    @{ Row = 86; F = $null; H = 67 },
    @{ Row = 88; F = $null; H = 67 },
    @{ Row = 89; F = $null; H = 67 },
    @{ Row = 91; F = $null; H = 67 },
    @{ Row = 96; F = "T3"; H = 45 },
    @{ Row = 99; F = "T20"; H = 60 },
    @{ Row = 102; F = "E4x10s"; H = 50 },
    @{ Row = 116; F = "T3"; H = 45 },
    @{ Row = 119; F = "T20"; H = 60 },
    @{ Row = 122; F = "E4x10s"; H = 40 }
)

foreach ($o in $overrides) {
    if ($o.F -ne $null) {
        $ws.Cells.Item($o.Row, 6).Value = $o.F
    }
    if ($o.H -ne $null) {
        $ws.Cells.Item($o.Row, 8).Value = $o.H
    }
}

# --- Apply a 2-decimal-place number format to the summary tables' ---
# --- Miles / Total Duration / Est. Overall Mileage columns (B:D)   ---
for ($r = 130; $r -le 141; $r++) {
    $ws.Range("B" + $r + ":D" + $r).NumberFormat = "0.00"
}
for ($r = 145; $r -le 160; $r++) {
    $ws.Range("B" + $r + ":D" + $r).NumberFormat = "0.00"
}
